$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 11:16"

# Row 15
$ws.Range("A15").Value = "Belgica"
$ws.Range("B15").Value = 2257
$ws.Range("C15").Value = 462
$ws.Range("D15").Value = 165
$ws.Range("E15").Value = 2055
$ws.Range("F15").Value = 130
$ws.Range("G15").Value = 16
$ws.Range("H15").Value = 37

# Row 16
$ws.Range("A16").Value = "Austria"
$ws.Range("B16").Value = 2203
$ws.Range("C16").Value = 24
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = 2188
$ws.Range("F16").Value = 13
$ws.Range("H16").Value = 6

# Row 17
$ws.Range("B17").Value = 1802
$ws.Range("C17").Value = 12
$ws.Range("E17").Value = 1794

# Row 20
$ws.Range("A20").Value = "Malasia"
$ws.Range("B20").Value = 1030
$ws.Range("C20").Value = 130
$ws.Range("D20").Value = 87
$ws.Range("E20").Value = 941
$ws.Range("F20").Value = 26
$ws.Range("H20").Value = 2

# Row 21
$ws.Range("A21").Value = "Japon"
$ws.Range("B21").Value = 963
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 215
$ws.Range("E21").Value = 715
$ws.Range("F21").Value = 49
$ws.Range("H21").Value = 33

# Row 27
$ws.Range("B27").Value = 705
$ws.Range("C27").Value = 28
$ws.Range("D27").Value = 15
$ws.Range("E27").Value = 690
$ws.Range("F27").Value = 10

# Row 88
$ws.Range("A88").Value = "Kazajistan"
$ws.Range("B88").Value = 49
$ws.Range("C88").Value = 5
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 49
$ws.Range("F88").Value = 0
$ws.Range("H88").Value = 0

# Row 89
$ws.Range("A89").Value = "Moldavia"
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 1
$ws.Range("E89").Value = 47
$ws.Range("F89").Value = 3
$ws.Range("H89").Value = 1

# Row 90
$ws.Range("A90").Value = "Lituania"
$ws.Range("B90").Value = 48
$ws.Range("F90").Value = 1
$ws.Range("H90").Value = 0

# Row 91
$ws.Range("A91").Value = "Oman"
$ws.Range("D91").Value = 13
$ws.Range("E91").Value = 35
$ws.Range("F91").Value = 0

# Row 92
$ws.Range("A92").Value = "Estado de Palestina"
$ws.Range("C92").Value = 1
$ws.Range("D92").Value = 17
$ws.Range("E92").Value = 31

# Row 94
$ws.Range("A94").Value = "Guadalupe"
$ws.Range("B94").Value = 45
$ws.Range("C94").Value = 12
$ws.Range("D94").Value = 0
$ws.Range("E94").Value = 45

# Row 95
$ws.Range("A95").Value = "Azerbaiyan"
$ws.Range("B95").Value = 44
$ws.Range("C95").Value = 0
$ws.Range("D95").Value = 7
$ws.Range("E95").Value = 36
$ws.Range("H95").Value = 1

# Row 96
$ws.Range("A96").Value = "Georgia"
$ws.Range("B96").Value = 43
$ws.Range("C96").Value = 3
$ws.Range("D96").Value = 1
$ws.Range("E96").Value = 42
$ws.Range("F96").Value = 1
$ws.Range("H96").Value = 0

# Row 97
$ws.Range("A97").Value = "Venezuela"
$ws.Range("B97").Value = 42
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("F97").Value = 0

# Row 98
$ws.Range("A98").Value = "Nueva Zelanda"
$ws.Range("B98").Value = 39
$ws.Range("C98").Value = 11
$ws.Range("E98").Value = 39

# Row 99
$ws.Range("A99").Value = "Tunez"
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 1
$ws.Range("E99").Value = 37
$ws.Range("F99").Value = 2
$ws.Range("H99").Value = 1

# Row 107
$ws.Range("D107").Value = 1
$ws.Range("E107").Value = 22

# Row 111
$ws.Range("A111").Value = "Banglades"
$ws.Range("B111").Value = 20
$ws.Range("C111").Value = 2
$ws.Range("D111").Value = 3
$ws.Range("E111").Value = 16
$ws.Range("H111").Value = 1

# Row 112
$ws.Range("A112").Value = "Consejo Danes para los Refugiados"
$ws.Range("C112").Value = 4
$ws.Range("D112").Value = 0
$ws.Range("E112").Value = 18
$ws.Range("H112").Value = 0

# Row 115
$ws.Range("A115").Value = "Ghana"
$ws.Range("E115").Value = 16
$ws.Range("H115").Value = 0

# Row 116
$ws.Range("A116").Value = "Cuba"
$ws.Range("C116").Value = 5
$ws.Range("D116").Value = 0
$ws.Range("E116").Value = 15

# Row 117
$ws.Range("A117").Value = "Jamaica"
$ws.Range("B117").Value = 16
$ws.Range("C117").Value = 1
$ws.Range("D117").Value = 2
$ws.Range("E117").Value = 13
$ws.Range("H117").Value = 1

# Row 118
$ws.Range("A118").Value = "Guayana Francesa"
$ws.Range("B118").Value = 15
$ws.Range("C118").Value = 0
$ws.Range("E118").Value = 15

# Row 119
$ws.Range("A119").Value = "Guam"
$ws.Range("B119").Value = 14
$ws.Range("C119").Value = 2
$ws.Range("E119").Value = 14

# Row 120
$ws.Range("A120").Value = "Montenegro"
$ws.Range("F120").Value = 0

# Row 121
$ws.Range("A121").Value = "Paraguay"
$ws.Range("F121").Value = 1

# Row 122
$ws.Range("A122").Value = "Maldivas"
$ws.Range("B122").Value = 13
$ws.Range("D122").Value = 0
$ws.Range("E122").Value = 13

# Row 123
$ws.Range("A123").Value = "Nigeria"
$ws.Range("B123").Value = 12
$ws.Range("D123").Value = 1

# Row 124
$ws.Range("A124").Value = "Monaco"
$ws.Range("C124").Value = 1

# Row 125
$ws.Range("A125").Value = "Ruanda"
$ws.Range("C125").Value = 0

# Row 126
$ws.Range("A126").Value = "Polinesia Francesa"
$ws.Range("C126").Value = 5
